# The commit swaps the contents of ppt/theme/theme1.xml ("Simple Light"
# color scheme) and ppt/theme/theme2.xml ("Default" color scheme).
# theme1.xml is the presentation's primary/active theme (used by the
# slide master and every slide), so the externally-visible effect of the
# swap is that the deck's color scheme changes from "Simple Light" to
# "Default". Apply that by rewriting the 12 theme colors on the active
# theme's color scheme to the "Default" scheme's values (dk1/lt1 are
# unchanged between the two schemes; dk2, lt2 and all six accents plus
# the two hyperlink colors change).

function HexToOle([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$tcs = $master.Theme.ThemeColorScheme

# Item order: 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2,
# 7 accent3, 8 accent4, 9 accent5, 10 accent6, 11 hlink, 12 folHlink
$tcs.Item(1).RGB = HexToOle "000000"
$tcs.Item(2).RGB = HexToOle "FFFFFF"
$tcs.Item(3).RGB = HexToOle "158158"
$tcs.Item(4).RGB = HexToOle "F3F3F3"
$tcs.Item(5).RGB = HexToOle "058DC7"
$tcs.Item(6).RGB = HexToOle "50B432"
$tcs.Item(7).RGB = HexToOle "ED561B"
$tcs.Item(8).RGB = HexToOle "EDEF00"
$tcs.Item(9).RGB = HexToOle "24CBE5"
$tcs.Item(10).RGB = HexToOle "64E572"
$tcs.Item(11).RGB = HexToOle "2200CC"
$tcs.Item(12).RGB = HexToOle "551A8B"
